$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update numeric values in column A
$ws.Range("A2").Value = 13
$ws.Range("A3").Value = 16
$ws.Range("A4").Value = 18

# Update the active selection to G12
$ws.Range("G12").Select()
